$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 46, shifting existing rows 46:102 down to 47:103
$ws.Rows("46:46").Insert()

$ws.Range("A46").Value = "JP26, JP27, JP28"
$ws.Range("C46").Value = "M05PTH (M05)"
$ws.Range("D46").Value = "1X05"
$ws.Range("E46").Value = "DK"
$ws.Range("F46").Value = "609-4303-ND"
$ws.Range("G46").Value = "68002-205HLF"
$ws.Range("H46").Value = 3
$ws.Range("I46").Value = 0.5
$ws.Range("J46").Formula = "=H46*I46"

$ws.Range("A16").Select()
$ws.Range("J46").Select()
